# Processing and dynamic mapping logic implemented
#
# - Rename existing blank "Sheet" to "StudentsMapping" and populate it with
#   each student's assigned organization per time slot.
# - Insert a new "OrganizationMapping" sheet (between StudentsMapping and the
#   new blank "Sheet") that is the inverse lookup: organization -> student
#   assigned to each time slot.
# - Re-add a blank "Sheet" at the end so the tab order/sheetId sequence is
#   preserved (StudentPreferences, StudentsMapping, OrganizationMapping, Sheet).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename the existing empty sheet to "StudentsMapping"
# ---------------------------------------------------------------------------
$studentsMapping = $wb.Worksheets.Item("Sheet")
$studentsMapping.Name = "StudentsMapping"

# ---------------------------------------------------------------------------
# 2. Insert "OrganizationMapping" right after "StudentsMapping"
# ---------------------------------------------------------------------------
$organizationMapping = $wb.Worksheets.Add($null, $studentsMapping)
$organizationMapping.Name = "OrganizationMapping"

# ---------------------------------------------------------------------------
# 3. Insert a fresh blank "Sheet" after "OrganizationMapping" (end of tabs)
# ---------------------------------------------------------------------------
$trailingSheet = $wb.Worksheets.Add($null, $organizationMapping)
$trailingSheet.Name = "Sheet"

# ---------------------------------------------------------------------------
# Populate StudentsMapping
# ---------------------------------------------------------------------------
$studentsMappingHeader = @("USC ID", "Student Name", "6:00 - 6:07pm", "6:10 - 6:17pm", "6:20 - 6:27pm", "6:30 - 6:37pm")
for ($col = 1; $col -le $studentsMappingHeader.Length; $col++) {
    $studentsMapping.Cells.Item(1, $col).Value = $studentsMappingHeader[$col - 1]
}

$studentsMappingRows = @(
    @(10001, "Oceana Hanner",       "Cedars-Sinai - Neurosciences", "Keck VIO - COBI",              "Optum CF - Patient XP",       "Providence Health Network"),
    @(10002, "Esther Choi",         "Keck VIO - COBI",              "Optum CF - Patient XP",         "St.Johns-PhysPartners ",      "Verdugo Hills Hospital"),
    @(10003, "Daniela Ahumada",     "City of Hope - CMO",           "CHLA - Anesthesia&CCM",         "Providence Health Network",   "Optum CF - Patient XP"),
    @(10004, "Fahima Gohil",        "Rancho Los Amigos NRC",        "Cedars-Sinai - Neurosciences",  "Emanate Health",              "City of Hope - CMO"),
    @(10005, "Julia Orozco",        "Emanate Health",               "Providence Health Network",     "Keck VIO - COBI",             "Kaiser PC - Consulting"),
    @(10006, "Emma Crusinberry",    "Kaiser PC - Consulting",       "St.Johns-PhysPartners ",        "SCAN Health Plan",            "Cedars-Sinai - Neurosciences"),
    @(10007, "Stanley Ibe",         "Optum CF - Patient XP",        "City of Hope - CMO",            "CHLA - Anesthesia&CCM",       "St.Johns-PhysPartners "),
    @(10008, "Raashi Subramanya",   "Providence Health Network",    "SCAN Health Plan",              "Cedars-Sinai - Neurosciences","Keck VIO - COBI")
)

$r = 2
foreach ($row in $studentsMappingRows) {
    for ($col = 1; $col -le $row.Length; $col++) {
        $studentsMapping.Cells.Item($r, $col).Value = $row[$col - 1]
    }
    $r++
}

# ---------------------------------------------------------------------------
# Populate OrganizationMapping
# ---------------------------------------------------------------------------
$organizationMappingHeader = @("Organization Code", "Organization Name", "6:00 - 6:07pm", "6:10 - 6:17pm", "6:20 - 6:27pm", "6:30 - 6:37pm")
for ($col = 1; $col -le $organizationMappingHeader.Length; $col++) {
    $organizationMapping.Cells.Item(1, $col).Value = $organizationMappingHeader[$col - 1]
}

$organizationMappingRows = @(
    @("C0", "Cedars-Sinai - Neurosciences",      10001,  10004,  10008,  10006),
    @("C1", "CHLA - Anesthesia&CCM",              $null,  10003,  10007,  $null),
    @("C2", "City of Hope - CMO",                 10003,  10007,   $null, 10004),
    @("E0", "Emanate Health",                     10005,   $null, 10004,  $null),
    @("K0", "Kaiser PC - Consulting",             10006,   $null,  $null, 10005),
    @("K1", "Keck IRM",                            $null,  $null,  $null, $null),
    @("K2", "Keck VIO - COBI",                    10002,  10001,  10005,  10008),
    @("O0", "Optum CF - Digi Transformation",      $null,  $null,  $null, $null),
    @("O1", "Optum CF - Patient XP",               10007,  10002,  10001,  10003),
    @("P0", "Providence Health Network",           $null,  $null,  $null, $null),
    @("P1", "Providence Health Network",           10008,  10005,  10003,  10001),
    @("R0", "Rancho Los Amigos NRC",               10004,   $null,  $null, $null),
    @("S0", "SCAN Health Plan",                    $null,  10008,  10006,  $null),
    @("S1", "St.Johns-PhysPartners ",               $null, 10006,  10002,  10007),
    @("T0", "Torrance Memorial",                   $null,  $null,  $null, $null),
    @("V0", "Verdugo Hills Hospital",              $null,  $null,   $null, 10002),
    @("W0", "West Hills Hospital",                 $null,  $null,  $null, $null)
)

$r = 2
foreach ($row in $organizationMappingRows) {
    for ($col = 1; $col -le $row.Length; $col++) {
        $val = $row[$col - 1]
        if ($null -ne $val) {
            $organizationMapping.Cells.Item($r, $col).Value = $val
        }
    }
    $r++
}
